$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin price / 1h-volume figures (cryptos list auto-update).
# Some new text values (e.g. "68.02") parse as plain numbers, but the
# source data keeps them as text (matching the original inlineStr cells),
# so those are round-tripped through a temporary Text number format and
# then restored to the Normal style to avoid leaving formatting residue.

$updates = @(
    @{ Cell = 'D2'; Value = '42.156.26'; KeepText = $true }
    @{ Cell = 'E2'; Value = '  +3.20%  '; KeepText = $true }
    @{ Cell = 'D3'; Value = '2.214.60'; KeepText = $true }
    @{ Cell = 'E4'; Value = '  -0.16%  '; KeepText = $true }
    @{ Cell = 'D5'; Value = '252.04'; KeepText = $false }
    @{ Cell = 'E5'; Value = '  +0.55%  '; KeepText = $true }
    @{ Cell = 'D6'; Value = '0.623'; KeepText = $false }
    @{ Cell = 'E6'; Value = '  +0.40%  '; KeepText = $true }
    @{ Cell = 'D7'; Value = '68.02'; KeepText = $false }
    @{ Cell = 'E7'; Value = '  +2.00%  '; KeepText = $true }
    @{ Cell = 'E8'; Value = '  -0.05%  '; KeepText = $true }
    @{ Cell = 'D9'; Value = '0.619'; KeepText = $false }
    @{ Cell = 'E9'; Value = '  +9.55%  '; KeepText = $true }
    @{ Cell = 'D10'; Value = '38.90'; KeepText = $false }
    @{ Cell = 'E10'; Value = '  +6.13%  '; KeepText = $true }
    @{ Cell = 'D11'; Value = '59.43'; KeepText = $false }
    @{ Cell = 'E11'; Value = '  +2.08%  '; KeepText = $true }
    @{ Cell = 'D12'; Value = '0.0938'; KeepText = $false }
    @{ Cell = 'E12'; Value = '  +1.60%  '; KeepText = $true }
    @{ Cell = 'D13'; Value = '7.03'; KeepText = $false }
    @{ Cell = 'E13'; Value = '  +1.60%  '; KeepText = $true }
    @{ Cell = 'E14'; Value = '  -0.31%  '; KeepText = $true }
    @{ Cell = 'D15'; Value = '2.546.81'; KeepText = $true }
    @{ Cell = 'E15'; Value = '  +1.93%  '; KeepText = $true }
    @{ Cell = 'D16'; Value = '0.869'; KeepText = $false }
    @{ Cell = 'E16'; Value = '  +1.59%  '; KeepText = $true }
    @{ Cell = 'D17'; Value = '14.51'; KeepText = $false }
    @{ Cell = 'E17'; Value = '  +1.45%  '; KeepText = $true }
    @{ Cell = 'D18'; Value = '2.210.79'; KeepText = $true }
    @{ Cell = 'E18'; Value = '  +1.85%  '; KeepText = $true }
    @{ Cell = 'D19'; Value = '41.997.36'; KeepText = $true }
    @{ Cell = 'E19'; Value = '  +2.97%  '; KeepText = $true }
    @{ Cell = 'D20'; Value = '0.0₃0962'; KeepText = $true }
    @{ Cell = 'E20'; Value = '  +2.62%  '; KeepText = $true }
    @{ Cell = 'D21'; Value = '72.31'; KeepText = $false }
    @{ Cell = 'E21'; Value = '  +1.43%  '; KeepText = $true }
    @{ Cell = 'D22'; Value = '6.14'; KeepText = $false }
    @{ Cell = 'E22'; Value = '  -0.22%  '; KeepText = $true }
    @{ Cell = 'D23'; Value = '231.21'; KeepText = $false }
    @{ Cell = 'E23'; Value = '  +0.41%  '; KeepText = $true }
    @{ Cell = 'E24'; Value = '  -1.17%  '; KeepText = $true }
    @{ Cell = 'D25'; Value = '3.88'; KeepText = $false }
    @{ Cell = 'E25'; Value = '  +0.56%  '; KeepText = $true }
    @{ Cell = 'E26'; Value = '  +0.16%  '; KeepText = $true }
    @{ Cell = 'D27'; Value = '11.18'; KeepText = $false }
    @{ Cell = 'E27'; Value = '  -3.98%  '; KeepText = $true }
    @{ Cell = 'E28'; Value = '  -2.42%  '; KeepText = $true }
    @{ Cell = 'D29'; Value = '3.71'; KeepText = $false }
    @{ Cell = 'E29'; Value = '  -0.63%  '; KeepText = $true }
    @{ Cell = 'D30'; Value = '2.20'; KeepText = $false }
    @{ Cell = 'E30'; Value = '  +1.19%  '; KeepText = $true }
    @{ Cell = 'D31'; Value = '166.69'; KeepText = $false }
    @{ Cell = 'E31'; Value = '  -1.44%  '; KeepText = $true }
    @{ Cell = 'D32'; Value = '20.43'; KeepText = $false }
    @{ Cell = 'E32'; Value = '  +0.06%  '; KeepText = $true }
    @{ Cell = 'D33'; Value = '5.95'; KeepText = $false }
    @{ Cell = 'E33'; Value = '  +11.26%  '; KeepText = $true }
    @{ Cell = 'E34'; Value = '  +3.90%  '; KeepText = $true }
    @{ Cell = 'D35'; Value = '0.0779'; KeepText = $false }
    @{ Cell = 'E35'; Value = '  +8.10%  '; KeepText = $true }
    @{ Cell = 'E36'; Value = '  +0.28%  '; KeepText = $true }
    @{ Cell = 'D37'; Value = '26.28'; KeepText = $false }
    @{ Cell = 'E37'; Value = '  +2.46%  '; KeepText = $true }
    @{ Cell = 'D38'; Value = '4.60'; KeepText = $false }
    @{ Cell = 'E38'; Value = '  +1.22%  '; KeepText = $true }
    @{ Cell = 'D39'; Value = '4.12'; KeepText = $false }
    @{ Cell = 'E39'; Value = '  +2.99%  '; KeepText = $true }
    @{ Cell = 'E40'; Value = '  +6.05%  '; KeepText = $true }
    @{ Cell = 'D41'; Value = '2.23'; KeepText = $false }
    @{ Cell = 'E41'; Value = '  +1.34%  '; KeepText = $true }
    @{ Cell = 'D42'; Value = '5.65'; KeepText = $false }
    @{ Cell = 'E42'; Value = '  +0.69%  '; KeepText = $true }
    @{ Cell = 'B43'; Value = 'Celestia'; KeepText = $true }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'; KeepText = $true }
    @{ Cell = 'D43'; Value = '12.03'; KeepText = $false }
    @{ Cell = 'E43'; Value = '  -1.57%  '; KeepText = $true }
    @{ Cell = 'B44'; Value = 'FTXToken'; KeepText = $true }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; KeepText = $true }
    @{ Cell = 'D44'; Value = '5.08'; KeepText = $false }
    @{ Cell = 'E44'; Value = '  +3.06%  '; KeepText = $true }
    @{ Cell = 'D45'; Value = '61.40'; KeepText = $false }
    @{ Cell = 'E45'; Value = '  -3.49%  '; KeepText = $true }
    @{ Cell = 'E46'; Value = '  -2.18%  '; KeepText = $true }
    @{ Cell = 'D47'; Value = '8.55'; KeepText = $false }
    @{ Cell = 'E47'; Value = '  +0.29%  '; KeepText = $true }
    @{ Cell = 'D48'; Value = '0.100'; KeepText = $false }
    @{ Cell = 'E48'; Value = '  -0.62%  '; KeepText = $true }
    @{ Cell = 'E49'; Value = '  -0.45%  '; KeepText = $true }
    @{ Cell = 'D50'; Value = '1.15'; KeepText = $false }
    @{ Cell = 'E50'; Value = '  +2.51%  '; KeepText = $true }
    @{ Cell = 'B51'; Value = 'SynthetixNetwork'; KeepText = $true }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'; KeepText = $true }
    @{ Cell = 'D51'; Value = '4.32'; KeepText = $false }
    @{ Cell = 'E51'; Value = '  +2.65%  '; KeepText = $true }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.KeepText) {
        # Plain string/unambiguous values: assigning .Value keeps them text.
        $rng.Value = $u.Value
    } else {
        # Numeric-looking text: force text storage, then restore the style
        # so no NumberFormat/style change lingers on the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    }
}
